$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "中超控股"
$ws.Range("A3").Value = "中文在线"
$ws.Range("B3").Value = "湖南白银"
$ws.Range("A4").Value = "巨力索具"
$ws.Range("B4").Value = "巨力索具"
$ws.Range("C4").Value = "巨力索具"
$ws.Range("A5").Value = "TCL中环"
$ws.Range("B5").Value = "TCL中环"
$ws.Range("C5").Value = "TCL中环"
$ws.Range("A6").Value = "湖南白银"
$ws.Range("B6").Value = "杉杉股份"
$ws.Range("C6").Value = "韩建河山"
$ws.Range("A7").Value = "杉杉股份"
$ws.Range("B7").Value = "中文在线"
$ws.Range("A8").Value = "中超控股"
$ws.Range("B8").Value = "利欧股份"
$ws.Range("C8").Value = "拓日新能"
$ws.Range("A9").Value = "利欧股份"
$ws.Range("B9").Value = "中超控股"
$ws.Range("A10").Value = "蓝色光标"
$ws.Range("B10").Value = "白银有色"
$ws.Range("C10").Value = "湖南白银"
$ws.Range("A11").Value = "拓日新能"
$ws.Range("B11").Value = "天孚通信"
$ws.Range("C11").Value = "完美世界"
$ws.Range("A12").Value = "杭电股份"
$ws.Range("B12").Value = "拓日新能"
$ws.Range("C12").Value = "横店影视"
$ws.Range("A13").Value = "数据港"
$ws.Range("B13").Value = "蓝色光标"
$ws.Range("C13").Value = "浙文互联"
$ws.Range("A14").Value = "白银有色"
$ws.Range("B14").Value = "三六零"
$ws.Range("C14").Value = "杉杉股份"
$ws.Range("A15").Value = "博纳影业"
$ws.Range("B15").Value = "杭电股份"
$ws.Range("C15").Value = "中文在线"
$ws.Range("A16").Value = "天孚通信"
$ws.Range("B16").Value = "神剑股份"
$ws.Range("C16").Value = "航天发展"
$ws.Range("A17").Value = "捷成股份"
$ws.Range("B17").Value = "捷成股份"
$ws.Range("C17").Value = "三六零"
$ws.Range("A18").Value = "三六零"
$ws.Range("B18").Value = "数据港"
$ws.Range("C18").Value = "白银有色"
$ws.Range("A19").Value = "航天发展"
$ws.Range("B19").Value = "航天发展"
$ws.Range("C19").Value = "神剑股份"
$ws.Range("A20").Value = "神剑股份"
$ws.Range("B20").Value = "浙文互联"
$ws.Range("C20").Value = "博纳影业"
$ws.Range("A21").Value = "浙文互联"
$ws.Range("B21").Value = "岩山科技"
$ws.Range("C21").Value = "蓝色光标"

